$d = $word.ActiveDocument

# Locate the "I don't know why this would be in italics..." paragraph (Block
# Text style). The two new paragraphs described by the diff must be inserted
# right after it:
#   1. "And this is also a proper way, with a different style" (FirstParagraph)
#   2. "This is called the Intense Quote style."              (BlockText / Intense Quote)

$anchorText = "I don" + [char]0x2019 + "t know why this would be in italics, but so it appears to be on my screen."

$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i).Range.Text.TrimEnd()
    if ($candidate -eq $anchorText) {
        $anchorIndex = $i
        break
    }
}
if ($anchorIndex -eq -1) {
    throw "Could not locate the anchor paragraph for the block-quote insertion."
}

$anchorPara = $d.Paragraphs.Item($anchorIndex)
$followingPara = $d.Paragraphs.Item($anchorIndex + 1)

# Insert the first new paragraph right before the paragraph that currently
# follows the anchor. It naturally inherits that paragraph's "FirstParagraph"
# style, so no explicit style assignment is required.
$followingPara.Range.InsertParagraphBefore()
$firstNew = $d.Paragraphs.Item($anchorIndex + 1)
$firstNew.Range.Text = "And this is also a proper way, with a different style"

# Insert the second new paragraph right after the first one. It inherits
# "FirstParagraph" from $firstNew, so it needs to be switched to the
# "Block Text" style (the Intense Quote mapping target).
$firstNew.Range.InsertParagraphAfter()
$secondNew = $d.Paragraphs.Item($anchorIndex + 2)
$secondNew.Range.Text = "This is called the Intense Quote style."
$secondNew.Style = "Block Text"
